# tests/schema/fixtures/invalid-data.xlsx
#
# 1) Rename the "Normals" sheet to "Normal records".
# 2) Turn on iterative calculation with a tighter max-change delta
#    (mirrors <calcPr iterateDelta="1E-4".../> in the saved workbook).
# 3) On that sheet, change header cell A1's text from "Id" to
#    "id with underscores" (adds a new shared string).
# 4) Update the sheet's saved selection from A4 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Normals")
$ws.Name = "Normal records"

$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

$ws.Range("A1").Value = "id with underscores"
$ws.Range("A2").Select() | Out-Null
